$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update/extend the houses data table (rows 10-32) ---
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 2100
$ws.Range("D10").Value = 230000
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = 2000
$ws.Range("D11").Value = 225000
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 1750
$ws.Range("D12").Value = 210000
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = 12
$ws.Range("C13").Value = 1850
$ws.Range("D13").Value = 225000
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = 13
$ws.Range("C14").Value = 1800
$ws.Range("D14").Value = 210000
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = 14
$ws.Range("C15").Value = 1700
$ws.Range("D15").Value = 200000
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = 15
$ws.Range("C16").Value = 1650
$ws.Range("D16").Value = 195000
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = 16
$ws.Range("C17").Value = 1900
$ws.Range("D17").Value = 240000
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = 17
$ws.Range("C18").Value = 1300
$ws.Range("D18").Value = 160000
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = 18
$ws.Range("C19").Value = 2500
$ws.Range("D19").Value = 325000
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = 19
$ws.Range("C20").Value = 1500
$ws.Range("D20").Value = 155000
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = 20
$ws.Range("C21").Value = 1400
$ws.Range("D21").Value = 170000
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = 21
$ws.Range("C22").Value = 2000
$ws.Range("D22").Value = 220000
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = 22
$ws.Range("C23").Value = 2135
$ws.Range("D23").Value = 245000
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = 23
$ws.Range("C24").Value = 2050
$ws.Range("D24").Value = 235000
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = 24
$ws.Range("C25").Value = 2250
$ws.Range("D25").Value = 275000
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = 25
$ws.Range("C26").Value = 1600
$ws.Range("D26").Value = 190000
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = 26
$ws.Range("C27").Value = 1950
$ws.Range("D27").Value = 230000
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = 27
$ws.Range("C28").Value = 2500
$ws.Range("D28").Value = 340000
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = 28
$ws.Range("C29").Value = 2200
$ws.Range("D29").Value = 300000
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = 29
$ws.Range("C30").Value = 1800
$ws.Range("D30").Value = 250000
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = 1350
$ws.Range("D31").Value = 145000
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 31
$ws.Range("C32").Value = 1250
$ws.Range("D32").Value = 145000

# --- Update the defined name range to cover the new extent ---
$houseName = $ws.Names.Item("houses")
$houseName.RefersTo = "=houses!`$A`$1:`$D`$32"

# --- Update the sheet view: scroll position and selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:B32").Select()
